$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44431
$ws.Range("M2").Value = 100

# Row 3
$ws.Range("D3").Value = 44405

# Row 4
$ws.Range("D4").Value = 44417
$ws.Range("M4").Value = 80

# Row 5
$ws.Range("D5").Value = 44418
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 1200
$ws.Range("O5").Value = 1200
$ws.Range("P5").Value = 1200
$ws.Range("S5").Value = 1200

# Row 6
$ws.Range("D6").Value = 44424
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 1200
$ws.Range("O6").Value = 1200
$ws.Range("P6").Value = 1200
$ws.Range("S6").Value = 1200

# Row 7
$ws.Range("D7").Value = 44343
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 1300
$ws.Range("O7").Value = 1300
$ws.Range("P7").Value = 1300
$ws.Range("S7").Value = 1300

# Row 8
$ws.Range("D8").Value = 44473
$ws.Range("M8").Value = 120

# Row 9
$ws.Range("D9").Value = 44435
$ws.Range("M9").Value = 130
$ws.Range("N9").Value = 1300
$ws.Range("O9").Value = 1300
$ws.Range("P9").Value = 1300
$ws.Range("S9").Value = 1300

# Row 10
$ws.Range("D10").Value = 44357
$ws.Range("M10").Value = 35
$ws.Range("N10").Value = 1000
$ws.Range("O10").Value = 1000
$ws.Range("P10").Value = 1000
$ws.Range("S10").Value = 1000

# Row 12
$ws.Range("D12").Value = 44432
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 1300
$ws.Range("O12").Value = 1300
$ws.Range("P12").Value = 1300
$ws.Range("S12").Value = 1300
